{"js": "// The SmartDocForm.docx is a data-bound form: three plain-text content\n// controls (structured document tags) show \"First Name\", \"Last Name\" and\n// \"Age\" values that are bound (via w:dataBinding) to the custom XML data\n// part {88E81A45-98C0-4D79-952A-E8203CE59AAC} (People/Person/...).\n//\n// Update the three bound values:\n//   First Name: Michael  -> Gisela\n//   Last Name:  Townsend -> Font\n//   Age:        33       -> 23\n\nconst contentControls = context.document.getContentControls();\ncontentControls.load(\"items/id\");\nawait context.sync();\n\n// Map each content control's id (from <w:sdt><w:sdtPr><w:id .../>) to its\n// new text value.\nconst newValuesById = {\n    \"458776938\": \"Gisela\", // First Name\n    \"930235383\": \"Font\",   // Last Name\n    \"57847020\": \"23\"       // Age\n};\n\nfor (const cc of contentControls.items) {\n    const newValue = newValuesById[cc.id];\n    if (newValue !== undefined) {\n        cc.insertText(newValue, Word.InsertLocation.replace);\n    }\n}\nawait context.sync();\n\n// Best effort: keep the bound custom XML data part in sync with the new\n// field values shown above (some hosts treat this part as read-only via\n// the JS API; ignore failures so the visible document edit still lands).\ntry {\n    const customXmlParts = context.document.customXmlParts;\n    customXmlParts.load(\"items/id\");\n    await context.sync();\n\n    const peopleCustomXmlPart = customXmlParts.items.find(\n        (p) => p.id === \"{88E81A45-98C0-4D79-952A-E8203CE59AAC}\"\n    );\n\n    if (peopleCustomXmlPart) {\n        const updatedXml =\n            '<People SDTemplate=\"People_Person.xml\">\\r\\n' +\n            \"  <Person>\\r\\n\" +\n            '    <fName Datatype=\"String\" Required=\"true\" DisplayName=\"First Name\">Gisela</fName>\\r\\n' +\n            '    <lName Datatype=\"String\" Required=\"true\" DisplayName=\"Last Name\">Font</lName>\\r\\n' +\n            '    <Age Datatype=\"Integer\" Required=\"false\" DisplayName=\"Age\">23</Age>\\r\\n' +\n            \"  </Person>\\r\\n\" +\n            \"</People>\";\n        peopleCustomXmlPart.setXml(updatedXml);\n        await context.sync();\n    }\n} catch (err) {\n    // Not supported in this host - the content control text above already\n    // reflects the requested change.\n}\n", "ps1": "# The SmartDocForm.docx is a data-bound form: three plain-text content\n# controls (structured document tags) show \"First Name\", \"Last Name\" and\n# \"Age\" values that are bound (via w:dataBinding) to the custom XML data\n# part {88E81A45-98C0-4D79-952A-E8203CE59AAC} (People/Person/...).\n#\n# Update the three bound values:\n#   First Name: Michael  -> Gisela\n#   Last Name:  Townsend -> Font\n#   Age:        33       -> 23\n\n$d = $word.ActiveDocument\n\nforeach ($cc in $d.ContentControls) {\n    switch ($cc.ID) {\n        \"458776938\" { $cc.Range.Text = \"Gisela\" }  # First Name\n        \"930235383\" { $cc.Range.Text = \"Font\" }    # Last Name\n        \"57847020\"  { $cc.Range.Text = \"23\" }      # Age\n    }\n}\n\n# Best effort: keep the bound custom XML data part in sync with the new\n# field values set above (some hosts treat this collection as read-only via\n# COM automation; ignore failures so the visible document edit still lands).\ntry {\n    $peopleStoreId = \"{88E81A45-98C0-4D79-952A-E8203CE59AAC}\"\n    $updatedXml = \"<People SDTemplate=`\"People_Person.xml`\">`r`n  <Person>`r`n    <fName Datatype=`\"String`\" Required=`\"true`\" DisplayName=`\"First Name`\">Gisela</fName>`r`n    <lName Datatype=`\"String`\" Required=`\"true`\" DisplayName=`\"Last Name`\">Font</lName>`r`n    <Age Datatype=`\"Integer`\" Required=`\"false`\" DisplayName=`\"Age`\">23</Age>`r`n  </Person>`r`n</People>\"\n\n    foreach ($part in $d.CustomXMLParts) {\n        if ($part.Id -eq $peopleStoreId) {\n            $part.XML = $updatedXml\n        }\n    }\n} catch {\n    # Not supported in this host - the content control text above already\n    # reflects the requested change.\n}\n"}
